# Change in setting of parameters in process within main and config file.
# Delete the "ReviewSheet_WorksheetName" / "Template" row from the Constants sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Locate the row to delete robustly by looking up the setting name in column A,
# in case row numbers ever drift; falls back to the known row (30).
$targetRow = 30
$found = $false
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    $cellValue = $ws.Cells.Item($r, 1).Value
    if ($cellValue -eq "ReviewSheet_WorksheetName") {
        $targetRow = $r
        $found = $true
        break
    }
}

$ws.Rows.Item($targetRow).Delete()

# Reflect the cursor position shown in the saved workbook after the edit.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A31").Select()
